# Auto-generated Excel COM-interop script to apply cryptos.xlsx (crypto price) update
# Commit: Updated cryptos list on Sat Jan 13 19:56:50 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '42.897.61'
$ws.Cells.Item(2, 5).Value = '  -1.51%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.578.67'
$ws.Cells.Item(3, 5).Value = '  -0.32%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '302.97'
$ws.Cells.Item(5, 5).Value = '  +0.61%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '96.48'
$ws.Cells.Item(6, 5).Value = '  +0.18%  '
$ws.Cells.Item(8, 5).Value = '  -0.08%  '
$ws.Cells.Item(9, 5).Value = '  -1.69%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '36.35'
$ws.Cells.Item(10, 5).Value = '  -0.57%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0812'
$ws.Cells.Item(11, 5).Value = '  -0.61%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '7.64'
$ws.Cells.Item(12, 5).Value = '  -1.31%  '
$ws.Cells.Item(13, 5).Value = '  +6.95%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '2.595.17'
$ws.Cells.Item(14, 5).Value = '  +0.46%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.887'
$ws.Cells.Item(15, 5).Value = '  -0.28%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '14.32'
$ws.Cells.Item(16, 5).Value = '  +0.01%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '42.979.77'
$ws.Cells.Item(17, 5).Value = '  -1.24%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '12.98'
$ws.Cells.Item(18, 5).Value = '  +5.27%  '
$ws.Cells.Item(19, 5).Value = '  +1.92%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.65'
$ws.Cells.Item(20, 5).Value = '  -0.12%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '71.99'
$ws.Cells.Item(21, 5).Value = '  -1.17%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '254.22'
$ws.Cells.Item(22, 5).Value = '  -4.03%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '2.95'
$ws.Cells.Item(23, 5).Value = '  +0.96%  '
$ws.Cells.Item(24, 5).Value = '  -3.04%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '29.01'
$ws.Cells.Item(25, 5).Value = '  -0.75%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.31'
$ws.Cells.Item(27, 5).Value = '  +0.26%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '37.54'
$ws.Cells.Item(28, 5).Value = '  -0.84%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.11'
$ws.Cells.Item(29, 5).Value = '  -2.30%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '6.06'
$ws.Cells.Item(30, 5).Value = '  -0.34%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '154.84'
$ws.Cells.Item(31, 5).Value = '  +1.67%  '
$ws.Cells.Item(32, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.44'
$ws.Cells.Item(32, 5).Value = '  -4.79%  '
$ws.Cells.Item(33, 2).Value = 'ARBITRUM'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '2.18'
$ws.Cells.Item(33, 5).Value = '  -1.10%  '
$ws.Cells.Item(34, 5).Value = '  -1.67%  '
$ws.Cells.Item(35, 5).Value = '  -0.77%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '18.22'
$ws.Cells.Item(36, 5).Value = '  +9.54%  '
$ws.Cells.Item(37, 5).Value = '  -3.04%  '
$ws.Cells.Item(38, 5).Value = '  -0.36%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '23.13'
$ws.Cells.Item(39, 5).Value = '  -4.56%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.15'
$ws.Cells.Item(40, 5).Value = '  +34.69%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.44'
$ws.Cells.Item(41, 5).Value = '  -4.64%  '
$ws.Cells.Item(42, 2).Value = 'RenderToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '3.91'
$ws.Cells.Item(42, 5).Value = '  +0.82%  '
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.0312'
$ws.Cells.Item(43, 5).Value = '  -0.70%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.088.52'
$ws.Cells.Item(44, 5).Value = '  +2.36%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.999'
$ws.Cells.Item(45, 5).Value = '  +0.19%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '9.22'
$ws.Cells.Item(46, 5).Value = '  +1.57%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '85.30'
$ws.Cells.Item(47, 5).Value = '  -2.97%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '76.75'
$ws.Cells.Item(48, 5).Value = '  +11.19%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '106.68'
$ws.Cells.Item(49, 5).Value = '  +1.01%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.827.42'
$ws.Cells.Item(50, 5).Value = '  -0.60%  '
$ws.Cells.Item(51, 5).Value = '  +2.49%  '
